$wb = $excel.ActiveWorkbook

# --- Step 1: insert the new '2022-Q1' sheet before '总计' ---
$zongjiSheet = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($zongjiSheet)
$q1.Name = "2022-Q1"

# --- Step 2: populate the '2022-Q1' sheet header row (B1:H1), styled like row 1 of sheet '2021-Q4' ---
$srcHeader = $wb.Worksheets.Item("2021-Q4").Range("B1:H1")
$srcHeader.Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# --- Step 3: populate the data rows (2..6) ---
$srcDataStyle = $wb.Worksheets.Item("2021-Q4").Range("A2:A6")
$srcDataStyle.Copy()
$q1.Range("A2:A6").PasteSpecial(-4122)

# row 2: 002379
$q1.Cells.Item(2, 1).Value = 0
$textRng = $q1.Range("B2:G2")
$textRng.NumberFormat = "@"
$q1.Cells.Item(2, 2).Value = "002379"
$q1.Cells.Item(2, 3).Value = "工银瑞信香港中小盘股票（QDII）人民币"
$q1.Cells.Item(2, 4).Value = "1.84"
$q1.Cells.Item(2, 5).Value = "86.48"
$q1.Cells.Item(2, 6).Value = "3.71"
$q1.Cells.Item(2, 7).Value = "0.0683"
$textRng.Style = "Normal"
$q1.Cells.Item(2, 8).Value = 7

# row 3: 002380
$q1.Cells.Item(3, 1).Value = 1
$textRng = $q1.Range("B3:G3")
$textRng.NumberFormat = "@"
$q1.Cells.Item(3, 2).Value = "002380"
$q1.Cells.Item(3, 3).Value = "工银瑞信香港中小盘股票（QDII）美元"
$q1.Cells.Item(3, 4).Value = "1.84"
$q1.Cells.Item(3, 5).Value = "86.48"
$q1.Cells.Item(3, 6).Value = "3.71"
$q1.Cells.Item(3, 7).Value = "0.0683"
$textRng.Style = "Normal"
$q1.Cells.Item(3, 8).Value = 7

# row 4: 012751
$q1.Cells.Item(4, 1).Value = 2
$textRng = $q1.Range("B4:G4")
$textRng.NumberFormat = "@"
$q1.Cells.Item(4, 2).Value = "012751"
$q1.Cells.Item(4, 3).Value = "建信纳斯达克100指数（QDII）A 美元现汇"
$q1.Cells.Item(4, 4).Value = "0.34"
$q1.Cells.Item(4, 5).Value = "88.02"
$q1.Cells.Item(4, 6).Value = "4.79"
$q1.Cells.Item(4, 7).Value = "0.0163"
$textRng.Style = "Normal"
$q1.Cells.Item(4, 8).Value = 5

# row 5: 012752
$q1.Cells.Item(5, 1).Value = 3
$textRng = $q1.Range("B5:G5")
$textRng.NumberFormat = "@"
$q1.Cells.Item(5, 2).Value = "012752"
$q1.Cells.Item(5, 3).Value = "建信纳斯达克100指数（QDII）C 人民币"
$q1.Cells.Item(5, 4).Value = "0.34"
$q1.Cells.Item(5, 5).Value = "88.02"
$q1.Cells.Item(5, 6).Value = "4.79"
$q1.Cells.Item(5, 7).Value = "0.0163"
$textRng.Style = "Normal"
$q1.Cells.Item(5, 8).Value = 5

# row 6: 012753
$q1.Cells.Item(6, 1).Value = 4
$textRng = $q1.Range("B6:G6")
$textRng.NumberFormat = "@"
$q1.Cells.Item(6, 2).Value = "012753"
$q1.Cells.Item(6, 3).Value = "建信纳斯达克100指数（QDII）C 美元现汇"
$q1.Cells.Item(6, 4).Value = "0.34"
$q1.Cells.Item(6, 5).Value = "88.02"
$q1.Cells.Item(6, 6).Value = "4.79"
$q1.Cells.Item(6, 7).Value = "0.0163"
$textRng.Style = "Normal"
$q1.Cells.Item(6, 8).Value = 5

# --- Step 4: rewrite the '总计' sheet with the new 2022-Q1 row prepended ---
$zongji = $wb.Worksheets.Item("总计")
# extend column-A styling (bold/border/centered, same as A2:A3) down to the new A4 row
$zongji.Range("A2").Copy()
$zongji.Range("A4").PasteSpecial(-4122)
$zongji.Cells.Item(2, 1).Value = 0
$zongji.Cells.Item(2, 2).Value = "2022-Q1"
$zongji.Cells.Item(2, 3).Value = 5
$zongji.Cells.Item(2, 4).Value = 0.19
$zongji.Cells.Item(3, 1).Value = 1
$zongji.Cells.Item(3, 2).Value = "2021-Q4"
$zongji.Cells.Item(3, 3).Value = 5
$zongji.Cells.Item(3, 4).Value = 1.28
$zongji.Cells.Item(4, 1).Value = 2
$zongji.Cells.Item(4, 2).Value = "2021-Q3"
$zongji.Cells.Item(4, 3).Value = 4
$zongji.Cells.Item(4, 4).Value = 0.02

Write-Output "edit complete"
